# Rebuttal letter edit ---------------------------------------------------
#
# "... have been moved to the supplementary results, as the highly
#  similar results underline robustness of our results."
# becomes
# "... have been moved to the supplementary results, but are highly
#  comparable to those obtained from AAL1."
#
# In the authored edit this happened as a normal Word editing session:
# the reviewer selected the clause "as the highly similar results
# underline robustness of our results" (leaving the leading ", " and
# the trailing "." alone) and retyped it. That is why the saved XML
# ends up with three runs instead of one, with Word's internal
# "last edit" bookmark (_GoBack), collapsed to zero length, sitting
# right after the freshly typed text and before the trailing period.

$d = $word.ActiveDocument

$old = "as the highly similar results underline robustness of our results"
$new = "but are highly comparable to those obtained from AAL1"

$range = $d.Content
$found = $range.Find.Execute($old, $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the target sentence to replace."
}

$editStart = $range.Start
$editEnd = $range.End

# Nudge the formatting of the character right after the clause (the
# trailing period) and back again. This forces that character to live
# in its own run boundary *before* anything else happens, the same way
# Word keeps the untouched trailing text in its own run once you start
# editing immediately in front of it.
$tail = $d.Range($editEnd, $editEnd + 1)
$tail.Font.Bold = $true
$tail.Font.Bold = $false

# Mark the two edges of the clause with temporary bookmarks so the
# clause becomes an isolated run on both sides, then retype its text.
$d.Bookmarks.Add("_TempEditStart", $d.Range($editStart, $editStart)) | Out-Null
$d.Bookmarks.Add("_TempEditEnd", $d.Range($editEnd, $editEnd)) | Out-Null

$clause = $d.Range($editStart, $editEnd)
$clause.Text = $new

$d.Bookmarks("_TempEditStart").Delete()
$d.Bookmarks("_TempEditEnd").Delete()

# Word keeps exactly one `_GoBack` bookmark in the document, always at
# the location of the most recent edit, collapsed to zero length right
# after the inserted text.
$clause.Collapse(0)
$d.Bookmarks.Add("_GoBack", $clause) | Out-Null
